$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 10417371
$ws.Range("I92").Value = 23809860
$ws.Range("J92").Value = 989.44446
$ws.Range("K92").Value = 23809860
$ws.Range("L92").Value = 989.44446
$ws.Range("M92").Value = -23808612
$ws.Range("N92").Value = -3485.44446

$ws.Range("H98").Value = 3373.75
$ws.Range("I98").Value = 2998.3333
$ws.Range("J98").Value = 4500
$ws.Range("K98").Value = 2998.3333
$ws.Range("L98").Value = 4500
$ws.Range("M98").Value = -1500.3333
$ws.Range("N98").Value = -7496

$ws.Range("H100").Value = 2787.2727
$ws.Range("I100").Value = 1365
$ws.Range("J100").Value = 4841.6665
$ws.Range("K100").Value = 1365
$ws.Range("L100").Value = 4841.6665
$ws.Range("M100").Value = -824
$ws.Range("N100").Value = -5923.6665

$ws.Range("H107").Value = 30303264
$ws.Range("I107").Value = 41666836
$ws.Range("K107").Value = 41666836
$ws.Range("M107").Value = -41664916

$ws.Range("H116").Value = 8097263.5
$ws.Range("I116").Value = 8335353.5
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 8335353.5
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -8331911.5
$ws.Range("N116").Value = -9084

$ws.Range("H122").Value = 3373.75
$ws.Range("I122").Value = 2998.3333
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 8994.999899999999
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -6544.999899999999
$ws.Range("N122").Value = -18400

$ws.Range("H132").Value = 3397.3845
$ws.Range("I132").Value = 3133.2896
$ws.Range("K132").Value = 9399.8688
$ws.Range("M132").Value = -6869.8688

$ws.Range("H137").Value = 1200.2639
$ws.Range("I137").Value = 939.7381
$ws.Range("J137").Value = 1565
$ws.Range("K137").Value = 2819.2143
$ws.Range("L137").Value = 4695
$ws.Range("M137").Value = -269.2143000000001
$ws.Range("N137").Value = -9795

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1026.9166
$ws.Range("I2").Value = 1025.8889
$ws.Range("J2").Value = 1030
$ws.Range("K2").Value = 1025.8889
$ws.Range("L2").Value = 1030
$ws.Range("M2").Value = -912.8888999999999
$ws.Range("N2").Value = -1256

$ws.Range("H32").Value = 8454.4
$ws.Range("I32").Value = 2578.4727
$ws.Range("J32").Value = 24613.2
$ws.Range("K32").Value = 2578.4727
$ws.Range("L32").Value = 24613.2
$ws.Range("M32").Value = -2291.4727
$ws.Range("N32").Value = -25187.2

$ws.Range("H45").Value = 1574
$ws.Range("I45").Value = 1579
$ws.Range("J45").Value = 1562.3334
$ws.Range("K45").Value = 1579
$ws.Range("L45").Value = 1562.3334
$ws.Range("M45").Value = -1202
$ws.Range("N45").Value = -2316.3334

$ws.Range("H110").Value = 5109.161
$ws.Range("I110").Value = 5347.5557
$ws.Range("J110").Value = 3500
$ws.Range("K110").Value = 5347.5557
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = -3302.5557
$ws.Range("N110").Value = -7590

$ws.Range("H116").Value = 1026.9166
$ws.Range("I116").Value = 1025.8889
$ws.Range("J116").Value = 1030
$ws.Range("K116").Value = 1025.8889
$ws.Range("L116").Value = 1030
$ws.Range("M116").Value = 1268.1111
$ws.Range("N116").Value = -5618

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1026.9166
$ws.Range("I3").Value = 1025.8889
$ws.Range("J3").Value = 1030
$ws.Range("K3").Value = 1025.8889
$ws.Range("L3").Value = 1030
$ws.Range("M3").Value = -911.8888999999999
$ws.Range("N3").Value = -1258

$ws.Range("H94").Value = 10123.087
$ws.Range("I94").Value = 1366.1666
$ws.Range("J94").Value = 41648
$ws.Range("K94").Value = 1366.1666
$ws.Range("L94").Value = 41648
$ws.Range("M94").Value = -915.1666
$ws.Range("N94").Value = -42550

$ws.Range("H107").Value = 2825.9412
$ws.Range("I107").Value = 2290.2
$ws.Range("J107").Value = 3591.2856
$ws.Range("K107").Value = 2290.2
$ws.Range("L107").Value = 3591.2856
$ws.Range("M107").Value = -370.1999999999998
$ws.Range("N107").Value = -7431.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 20834910
$ws.Range("I16").Value = 45456344
$ws.Range("J16").Value = 1387.2307
$ws.Range("K16").Value = 45456344
$ws.Range("L16").Value = 1387.2307
$ws.Range("M16").Value = -45456057
$ws.Range("N16").Value = -1961.2307

$ws.Range("H31").Value = 7548938.5
$ws.Range("I31").Value = 9525078
$ws.Range("J31").Value = 6252097
$ws.Range("K31").Value = 9525078
$ws.Range("L31").Value = 6252097
$ws.Range("M31").Value = -9524783
$ws.Range("N31").Value = -6252687

$ws.Range("H34").Value = 7548938.5
$ws.Range("I34").Value = 9525078
$ws.Range("J34").Value = 6252097
$ws.Range("K34").Value = 9525078
$ws.Range("L34").Value = 6252097
$ws.Range("M34").Value = -9524876
$ws.Range("N34").Value = -6252501

$ws.Range("H58").Value = 1572.9333
$ws.Range("I58").Value = 759.8
$ws.Range("J58").Value = 1979.5
$ws.Range("K58").Value = 759.8
$ws.Range("L58").Value = 1979.5
$ws.Range("M58").Value = -556.8
$ws.Range("N58").Value = -2385.5

$ws.Range("H113").Value = 20834910
$ws.Range("I113").Value = 45456344
$ws.Range("J113").Value = 1387.2307
$ws.Range("K113").Value = 45456344
$ws.Range("L113").Value = 1387.2307
$ws.Range("M113").Value = -45454174
$ws.Range("N113").Value = -5727.2307

$ws.Range("H122").Value = 4168288.5
$ws.Range("J122").Value = 1730.2
$ws.Range("L122").Value = 5190.6
$ws.Range("N122").Value = -10090.6

$ws.Range("H132").Value = 3089.6
$ws.Range("I132").Value = 974.75
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 2924.25
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -394.25
$ws.Range("N132").Value = -18558.5

$ws.Range("H136").Value = 1572.9333
$ws.Range("I136").Value = 759.8
$ws.Range("J136").Value = 1979.5
$ws.Range("K136").Value = 2279.4
$ws.Range("L136").Value = 5938.5
$ws.Range("M136").Value = 270.6000000000004
$ws.Range("N136").Value = -11038.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 7936887.5
$ws.Range("I107").Value = 20833646
$ws.Range("J107").Value = 421.07693
$ws.Range("K107").Value = 20833646
$ws.Range("L107").Value = 421.07693
$ws.Range("M107").Value = -20831726
$ws.Range("N107").Value = -4261.07693

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1702.2941
$ws.Range("I93").Value = 1517.1818
$ws.Range("J93").Value = 2041.6666
$ws.Range("K93").Value = 1517.1818
$ws.Range("L93").Value = 2041.6666
$ws.Range("M93").Value = -269.1818000000001
$ws.Range("N93").Value = -4537.6666

$ws.Range("H122").Value = 2702.8572
$ws.Range("I122").Value = 2640
$ws.Range("J122").Value = 2737.7778
$ws.Range("K122").Value = 7920
$ws.Range("L122").Value = 8213.3334
$ws.Range("M122").Value = -5470
$ws.Range("N122").Value = -13113.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 52632850
$ws.Range("I81").Value = 76924050
$ws.Range("J81").Value = 1915.6666
$ws.Range("K81").Value = 153848100
$ws.Range("L81").Value = 3831.3332
$ws.Range("M81").Value = -153847039
$ws.Range("N81").Value = -5953.3332

$ws.Range("H84").Value = 52632850
$ws.Range("I84").Value = 76924050
$ws.Range("J84").Value = 1915.6666
$ws.Range("K84").Value = 769240500
$ws.Range("L84").Value = 19156.666
$ws.Range("M84").Value = -769235196
$ws.Range("N84").Value = -29764.666

Write-Output "All changes applied"
